$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Files:" column (J2:J4) with the new cgap-style accession values
$ws.Range("J2").Value = "cgap:NA12879_sample_S1_R1, cgap:NA12879_sample_S1_R2"
$ws.Range("J3").Value = "cgap:NA12878_sample_S1_R1, cgap:NA12878_sample_S1_R2"
$ws.Range("J4").Value = "cgap:NA12877_sample_S1_R1, cgap:NA12877_sample_S1_R2"

# Move the selection back to A2 (single cell) instead of J2:J4
$ws.Range("A2").Select()
